$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "63.303.77"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.445.20"
$ws.Range("E3").Value = "  +0.00%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "571.24"
$ws.Range("E5").Value = "  +0.69%  "
Set-TextValue $ws.Range("D6") "146.87"
$ws.Range("E6").Value = "  +0.55%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.539"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "2.440.88"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("E12").Value = "  -1.30%  "
Set-TextValue $ws.Range("D13") "0.354"
$ws.Range("E13").Value = "  -0.30%  "
Set-TextValue $ws.Range("D14") "27.07"
$ws.Range("E14").Value = "  +0.08%  "
Set-TextValue $ws.Range("D15") "0.0000178"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "63.052.93"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "2.431.82"
$ws.Range("E18").Value = "  +0.03%  "
Set-TextValue $ws.Range("D19") "11.30"
$ws.Range("E19").Value = "  +0.09%  "
Set-TextValue $ws.Range("D20") "7.35"
$ws.Range("E20").Value = "  +5.75%  "
Set-TextValue $ws.Range("D21") "327.77"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  +0.15%  "
Set-TextValue $ws.Range("D23") "2.07"
$ws.Range("E23").Value = "  +11.41%  "
$ws.Range("E24").Value = "  +4.36%  "
Set-TextValue $ws.Range("D25") "65.31"
$ws.Range("E25").Value = "  -3.07%  "
Set-TextValue $ws.Range("D26") "616.72"
$ws.Range("E26").Value = "  +4.06%  "
Set-TextValue $ws.Range("D27") "8.82"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.559.55"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D30") "1.50"
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  +0.16%  "
Set-TextValue $ws.Range("D32") "8.28"
$ws.Range("E32").Value = "  -2.20%  "
Set-TextValue $ws.Range("D33") "0.141"
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("E34").Value = "  +0.94%  "
Set-TextValue $ws.Range("D35") "5.20"
$ws.Range("E35").Value = "  +6.65%  "
Set-TextValue $ws.Range("D36") "1.52"
$ws.Range("E36").Value = "  -1.56%  "
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.18%  "
Set-TextValue $ws.Range("D38") "0.380"
$ws.Range("E38").Value = "  -0.90%  "
Set-TextValue $ws.Range("D39") "5.42"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D41") "2.69"
$ws.Range("E41").Value = "  +10.40%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D42") "145.57"
$ws.Range("E42").Value = "  -1.68%  "
Set-TextValue $ws.Range("D43") "1.79"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("E44").Value = "  -0.50%  "
Set-TextValue $ws.Range("D45") "148.80"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  +1.91%  "
Set-TextValue $ws.Range("D47") "21.22"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("E48").Value = "  -0.12%  "
Set-TextValue $ws.Range("D49") "0.600"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +0.36%  "
Set-TextValue $ws.Range("D51") "0.0915"
$ws.Range("E51").Value = "  -1.20%  "
